$d = $word.ActiveDocument

# Locate the "Difficulty: Medium" paragraph (a standalone bulleted line near
# the top of the article) and remove the whole paragraph, including its
# paragraph mark, so the following Subtitle paragraph moves up to take its
# place.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Difficulty:*") {
        $para.Range.Delete()
        $found = $true
        break
    }
}

if ($found) {
    # Word stamps the `_GoBack` bookmark at the location of the most recent
    # edit; re-create that behavior at the start of the paragraph that now
    # follows the deleted bullet (the article subtitle run).
    $afterPara = $d.Paragraphs.Item($i)
    $gobackRange = $d.Range($afterPara.Range.Start, $afterPara.Range.Start)

    foreach ($bm in $d.Bookmarks) {
        if ($bm.Name -eq "_GoBack") {
            $bm.Delete()
        }
    }

    $d.Bookmarks.Add("_GoBack", $gobackRange)
}
